$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 21387
$ws1.Range("F3").Value = 3268
$ws1.Range("F4").Value = 849
$ws1.Range("F8").Value = 300
$ws1.Range("F15").Value = 35
$ws1.Range("F16").Value = 451
$ws1.Range("F17").Value = 160
$ws1.Range("F19").Value = 30

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F4").Value = 724
$ws3.Range("F5").Value = 1723
$ws3.Range("F6").Value = 78

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 724
$ws4.Range("F5").Value = 1723
$ws4.Range("F6").Value = 21388
$ws4.Range("F7").Value = 3268
$ws4.Range("F8").Value = 849
$ws4.Range("F10").Value = 78
$ws4.Range("F14").Value = 300
$ws4.Range("F29").Value = 35
$ws4.Range("F30").Value = 451
$ws4.Range("F32").Value = 160
$ws4.Range("F36").Value = 30
